$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 538.1254780066262
$ws.Range("D2").Value = 120.1197662093448
$ws.Range("F2").Value = 446
$ws.Range("G2").Value = 504
$ws.Range("H2").Value = 593
$ws.Range("C3").Value = 43.19081207395062
$ws.Range("D3").Value = 5.001223665698476
$ws.Range("E3").Value = 28.54
$ws.Range("F3").Value = 39.77
$ws.Range("G3").Value = 43.1
$ws.Range("H3").Value = 46.59
$ws.Range("C4").Value = 1.564071704035329
$ws.Range("D4").Value = 2.784652578677638
$ws.Range("F4").Value = 0.57
$ws.Range("G4").Value = 1.16
$ws.Range("H4").Value = 2.1
$ws.Range("C5").Value = 322.3611026898137
$ws.Range("D5").Value = 10.44548591421454
$ws.Range("F5").Value = 315.98
$ws.Range("G5").Value = 324.38
$ws.Range("H5").Value = 330.69
$ws.Range("C6").Value = 22.11743291686436
$ws.Range("D6").Value = 1.942015821583213
$ws.Range("F6").Value = 20.91
$ws.Range("G6").Value = 21.96
$ws.Range("H6").Value = 22.87
$ws.Range("C7").Value = -76.83377129486571
$ws.Range("D7").Value = 23.67980839532934
$ws.Range("F7").Value = -94
$ws.Range("C8").Value = 7.519434189096432
$ws.Range("D8").Value = 6.887967962223891
$ws.Range("C9").Value = 9.32322169369457
$ws.Range("D9").Value = 1.688405586861373
$ws.Range("C10").Value = 867.8304302161014
$ws.Range("D10").Value = 0.4613076307960154
$ws.Range("C11").Value = 0.5571491068585587
$ws.Range("D11").Value = 0.5908570403400689
$ws.Range("C12").Value = 22.71709509848052
$ws.Range("D12").Value = 12.2911372628886
$ws.Range("C13").Value = 0.6725933405515883
$ws.Range("D13").Value = 0.7488690246624591
$ws.Range("C14").Value = 1.829280161559143
$ws.Range("D14").Value = 1.66849122640517
$ws.Range("C15").Value = 94.23377129486568
$ws.Range("D15").Value = 23.67980839532417
$ws.Range("H15").Value = 111.4
$ws.Range("C16").Value = -85.94257698204969
$ws.Range("D16").Value = 21.33025101082983
$ws.Range("E16").Value = -137.0738221927363
$ws.Range("F16").Value = -103.3377954106368
$ws.Range("G16").Value = -84.69305820175224
$ws.Range("H16").Value = -66.33195619988427
$ws.Range("C17").Value = -78.42314279295324
$ws.Range("D17").Value = 25.97442681196299
$ws.Range("F17").Value = -93.93380807687734
$ws.Range("G17").Value = -74.21238401914255
$ws.Range("H17").Value = -55.43249407632486
